$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row (row 3 is currently blank) ------------------------------
$ws.Range("C3").Value2 = "Included"
$ws.Range("J3").Value2 = "Removed"
$ws.Range("C3:J3").Font.Bold = $true

# --- 2. Update the sequoia-pgp entry in place (row 19) ---------------------
# keeps its existing hyperlink relationship (target stays https://github.com/sequoia-pgp/)
$ws.Range("C19").Value2 = "https://github.com/sequoia-pgp/sequoia-git"

# --- 3. Move four URLs out of column C into column J ("Removed") -----------
# capture happens implicitly: we already know their text/targets from the
# original sheet, so just add them directly as new hyperlinks in column J.
$ws.Hyperlinks.Add($ws.Range("J4"), "https://github.com/tc39/", "", "", "https://github.com/tc39/")
$ws.Hyperlinks.Add($ws.Range("J5"), "https://github.com/rustls/", "", "", "https://github.com/rustls/ ")
$ws.Hyperlinks.Add($ws.Range("J6"), "https://github.com/php/", "", "", "https://github.com/php/")
$ws.Hyperlinks.Add($ws.Range("J7"), "https://github.com/openjs-foundation/", "", "", "https://github.com/openjs-foundation/ ")
$ws.Range("J4:J7").Style = "Hyperlink"

# --- 4. Remove the now-relocated / deleted rows from column C --------------
# delete bottom-to-top so the row numbers referenced above stay valid
$ws.Rows("35:35").Delete()
$ws.Rows("27:27").Delete()
$ws.Rows("17:17").Delete()
$ws.Rows("12:12").Delete()
$ws.Rows("11:11").Delete()
